$wb = $excel.ActiveWorkbook

# --- Product Backlog sheet: mark the last 4 backlog items as "Done" ---
$pb = $wb.Worksheets.Item("Product Backlog")

$doneRange = $pb.Range("I20:I23")
$doneRange.Value = "Done"
# Re-apply the same look-and-feel used by the other "Done" cells (e.g. I19):
# green "Good" text/fill, centered, thin border already present on these cells.
$doneRange.Font.Color = 24832        # RGB(0,97,0)  == FF006100 (OLE BGR order)
$doneRange.Interior.Color = 13561798 # RGB(198,239,206) == FFC6EFCE (OLE BGR order)
$doneRange.HorizontalAlignment = -4108 # xlCenter

# The old "In Progress" cell style is no longer used anywhere in the workbook,
# so drop the now-orphaned "Neutral" cell style definition.
$wb.Styles.Item("Neutral").Delete()

# --- Sprint sheet: fill in day-3 ("Hari 3") actuals for Sprint 6's burndown log ---
$sp = $wb.Worksheets.Item("Sprint")

$sp.Range("F115").Value = 0
$sp.Range("F116").Value = 0
$sp.Range("F117").Value = 0
$sp.Range("F118").Value = 0
$sp.Range("F119").Value = 0
$sp.Range("F120").Value = 0
$sp.Range("F121").Value = 0
$sp.Range("F122").Value = 0
$sp.Range("F123").Value = 1
$sp.Range("F124").Value = 1
$sp.Range("F125").Value = 1

# Burndown total for day 3 (mirrors the existing D128/E128 pattern).
$sp.Range("F128").Formula = "=E128-(SUM(F115:F125))"
